$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "EMP_1"
$ws.Range("C2").Value = "ryann@mail.com"
$ws.Range("E2").Value = "BSIT"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "09123456789"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1000"
$ws.Range("H2").Value = "Sa bahay"

# Row 3
$ws.Range("A3").Value = "EMP_2"
$ws.Range("B3").Value = "Malabanan, RySes"
$ws.Range("C3").Value = "ryses@email.com"
$ws.Range("D3").Value = "male"
$ws.Range("E3").Value = "BSIT"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "09123456789"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "160"
$ws.Range("H3").Value = "Sa bahay"

# Row 4
$ws.Range("A4").Value = "EMP_3"
$ws.Range("B4").Value = "Malabanan, Kim"
$ws.Range("C4").Value = "kim@email.com"
$ws.Range("E4").Value = "BSIT"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "100"
$ws.Range("H4").Value = "Sample"
